# Update Wilke validation values (column C, rows 2-9) on the "Ark1" sheet.
# Column D (AMS %) holds formulas (=Cn/C5) that recalculate automatically
# once the inputs change, and the bar chart reads those same D2:D9 cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("C2").Value = 102.37860000000001
$ws.Range("C3").Value = 304.15249999999997
$ws.Range("C4").Value = 522.53359999999998
$ws.Range("C5").Value = 556.05790000000002
$ws.Range("C6").Value = 1273.777
$ws.Range("C7").Value = 1183.127
$ws.Range("C8").Value = 2636.0650000000001
$ws.Range("C9").Value = 2123.0590000000002

# Matches the author's last-saved selection in the updated workbook.
$ws.Range("G26").Select()

$excel.CalculateFull()
